# Rename the "Collection" sheet tab to "CRF" (commit: "rename Collection to CRF in tabs")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "CRF_SC"
